# Generate Report for Handoff
# Update the "Latest Handoff/Generate" timestamps for the row corresponding to
# ac892881-cb09-44a5-a814-f764cdaf6522.md on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the last data row (row 7)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-29 22:44:00"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the last data row (row 7)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-29 22:43:56"

# de-de sheet: "Latest Handoff Datetime" column (H) for the last data row (row 7)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-29 22:44:00"
